# Apply the "output generated at 456a3b4" data refresh to both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, which carry
# duplicate copies of the same event table.
#
# Per-row changes (想去人数 = F, 最低票价 = G, 名称 = C):
#   F3:  1446 -> 1450
#   F9:  197  -> 198
#   F10: 142  -> 143
#   F12: 4784 -> 4800
#   F14: 7069 -> 7075
#   C15: 赣州·马娘only -> 赣州·赛马娘only
#   F18: 583  -> 585
#   F20: 9    -> 12
#   F21: 4190 -> 4196 ; G21: 64 -> 54
#   F22: 1320 -> 1359
#   F23: 84   -> 85
#   F24: 78   -> 79
#   F25: 2770 -> 2772
#   F29: 396  -> 399
#   F30: 394  -> 398
#   F31: 418  -> 419
#   F32: 248  -> 250
#   F33: 56   -> 57
#   F34: 1655 -> 1658
#   F35: 1069 -> 1073
#   F36: 75   -> 76
#   F37: 923  -> 957
#   F38: 90   -> 91
#   F39: 559  -> 560
#   F41: 501  -> 502
#   F42: 10   -> 11
#   F43: 24   -> 26
#   F45: sheet1 852 -> 968 ; sheet4 853 -> 969 (sheets differ by 1 already)
#   F47: 26   -> 27

$wb = $excel.ActiveWorkbook

# Numeric "想去人数" (F column) updates that are identical across both
# data sheets.
$fUpdates = @{
    3  = 1450
    9  = 198
    10 = 143
    12 = 4800
    14 = 7075
    18 = 585
    20 = 12
    21 = 4196
    22 = 1359
    23 = 85
    24 = 79
    25 = 2772
    29 = 399
    30 = 398
    31 = 419
    32 = 250
    33 = 57
    34 = 1658
    35 = 1073
    36 = 76
    37 = 957
    38 = 91
    39 = 560
    41 = 502
    42 = 11
    43 = 26
    47 = 27
}

# "最低票价" (G column) update(s).
$gUpdates = @{
    21 = 54
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
    }

    # C15 name correction: 赣州·马娘only -> 赣州·赛马娘only
    $ws.Cells.Item(15, 3).Value = "赣州·赛马娘only"
}

# Row 45's F value differs by 1 between the two sheets both before and
# after the refresh, so set it per-sheet explicitly (overrides nothing
# above since row 45 wasn't in $fUpdates).
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(45, 6).Value = 968

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(45, 6).Value = 969
